$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EventMonitoring")

$ws.Range("E2").Value = 'Check RAMP to prevent FOD.
1. Quan sát việc tuân thủ kiểm tra FOD bãi đậu trước khi tàu đáp và sau khi đưa tàu đi của NVKT.
2. Croscheck tại các bãi có bảo dưỡng lớn xem team bải dưỡng có tuân thủ quy trình ngăn ngừa FOD hay không?'

$ws.Range("E3").Value = 'Check fuel drainage according to the checklist.Random check các tàu thực hiện WKLY có được drain fuel theo WKLY checklist hay không?'

$ws.Range("E4").Value = 'Ensure OWP feedback is fully recorded.Có thể kiểm tra OWP feedback của ngày hôm trước, với các tiêu chí:
        a. Các WO từ chối phải có lí do rõ ràng và được sự xác nhận từ MOC.
        b. Thông tin daily check phải được điền đầy đủ.
        c. Các WO được thực hiện phải ghi nhận số chứng chỉ của NVKT rõ ràng.'

$ws.Range("E7").Value = 'Check for damage, leaks, and defects in the cargo.
a. Kiểm tra các dấu hiệu mục của buồng hàng - đặc biệt là khu vực quanh mép buồng hàng và cửa buồng hàng.
        b. Kiểm tra tình trạng TDP.
        c. Kiểm tra các tấm linning buồng hàng.'

$ws.Range("E8").Value = 'Ensure engine wash chemicals are used correctly.Chi tiết kiểm tra:
        a. Số lượng hóa chất.
        b. Số lượng lần rửa'

$ws.Range("E9").Value = 'Conduct a general inspection of the cabin, Lavatory, and Galley. Tăng cường random check áo phao'

$ws.Range("E12").Value = 'Check for engine failures, leaks, and loose bolts. Tăng cường kiểm tra vị trí Anti-ice access panel'

$ws.Range("E20").Value = 'Verify compliance with PTS in operations. PTS là chương trình yêu cầu các đơn vị phải hoàn thành phần việc của mình trong khoảng thời gian CỐ ĐỊNH và TỐI ƯU NHẤT đã được thống nhất bằng văn bản. Nhằm duy trì tổng thời gian dừng/nghỉ giữa các chuyến bay TỐI ƯU.'

$ws.Range("E5").Select()
